$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shared strings: "N группа" -> Roman numeral equivalents
$ws.Range("B2").Value = "I группа"
$ws.Range("B3").Value = "II группа"
$ws.Range("B4").Value = "III группа"
$ws.Range("B5").Value = "IV группа"
$ws.Range("B6").Value = "V группа"

# Update C and D columns for rows 2-6
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 3).Value = 2
    $ws.Cells.Item($r, 4).Value = 2958100
}
